$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 81 with the corrected Start Time and the missing Stop Time,
# Interruption and Activity data
$ws.Range("B81").Value = 0.59305555555555556
$ws.Range("C81").Value = 0.70624999999999993
$ws.Range("D81").Value = 30
$ws.Range("E81").Formula = '=IF(AND(NOT(ISBLANK(B81)),NOT(ISBLANK(C81))),(C81-B81)*24-D81/60,"")'
$ws.Range("F81").Value = "Coding"

$excel.CalculateFull()

# Update the selection to reflect where the user left off editing
$ws.Range("B82").Select()
